# Apply cryptos list update (diff-driven edits)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.894.03"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.924.21"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "240.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("E6").Value = "  +0.11%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4913"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2977"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06782"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "1.913.52"
$ws.Range("E10").Value = "  +1.23%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "17.10"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07303"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.81%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.180"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.11%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "89.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6736"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "30.867.92"
$ws.Range("E16").Value = "  +0.60%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008000"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.57"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.50%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "2.171.08"
$ws.Range("E20").Value = "  +1.66%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.191"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +7.22%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "206.48"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +8.33%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.311"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.704"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.36%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "159.15"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.13"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.994"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.13%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.428"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.367"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09195"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("E32").Value = "  +1.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05208"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7564"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.127"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.72%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.730"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.25%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01862"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.739"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.30%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.9281"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.100"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.06%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.4531"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.10%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "108.30"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.932"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "70.88"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +22.20%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.012"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.16%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.1401"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.30%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "7.708"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.104"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.83%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "35.55"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.05%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.4108"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.92%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05957"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
